$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J (pushes TL Borc/TL Alacak/TL Bakiye/B-A and the
# Toplam: cell one column to the right) and label it "Kur".
$ws.Columns("J").Insert()
$ws.Range("J4").Value = "Kur"

# Restore the active selection to the cell the author ended up on.
$ws.Range("I18").Select()
